$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: copy formatting (style) from column BA (rows 1-80) to the new BB:BD columns
$fmtSrc = $ws.Range("BA1:BA80")
$fmtDst = $ws.Range("BB1:BD80")
$fmtSrc.Copy($fmtDst)

# Step 2: write the actual values for BB, BC, BD (31/12/2023, 31/03/2024, 30/06/2024 columns)
# Build a single 80x3 array (rows 1..80, columns BB..BD) and assign in one shot.
$data = New-Object 'object[,]' 80,3
$data[0,0] = "31/12/2023"
$data[0,1] = "31/03/2024"
$data[0,2] = "30/06/2024"
$data[1,0] = 4509324.8
$data[1,1] = 4565757.952
$data[1,2] = 4330566.144
$data[2,0] = 1618706.048
$data[2,1] = 1760907.008
$data[2,2] = 1570237.952
$data[3,0] = 394588.992
$data[3,1] = 497505.984
$data[3,2] = 395671.008
$data[4,0] = 592179.008
$data[4,1] = 629643.008
$data[4,2] = 399620.992
$data[5,0] = 242992.992
$data[5,1] = 237071.008
$data[5,2] = 296801.984
$data[6,0] = 0
$data[6,1] = 0
$data[6,2] = 0
$data[7,0] = 0
$data[7,1] = 0
$data[7,2] = 0
$data[8,0] = 0
$data[8,1] = 0
$data[8,2] = 0
$data[9,0] = 0
$data[9,1] = 0
$data[9,2] = 0
$data[10,0] = 388944.992
$data[10,1] = 396687.008
$data[10,2] = 478144
$data[11,0] = 297403.008
$data[11,1] = 283140
$data[11,2] = 294417.984
$data[12,0] = 0
$data[12,1] = 0
$data[12,2] = 0
$data[13,0] = 0
$data[13,1] = 0
$data[13,2] = 0
$data[14,0] = 0
$data[14,1] = 0
$data[14,2] = 0
$data[15,0] = 0
$data[15,1] = 0
$data[15,2] = 0
$data[16,0] = 0
$data[16,1] = 0
$data[16,2] = 0
$data[17,0] = 0
$data[17,1] = 0
$data[17,2] = 0
$data[18,0] = 204066
$data[18,1] = 204032
$data[18,2] = 218492.992
$data[19,0] = 0
$data[19,1] = 0
$data[19,2] = 0
$data[20,0] = 0
$data[20,1] = 0
$data[20,2] = 0
$data[21,0] = 262
$data[21,1] = 262
$data[21,2] = 262
$data[22,0] = 38152
$data[22,1] = 31494
$data[22,2] = 40495
$data[23,0] = 2554801.92
$data[23,1] = 2489955.072
$data[23,2] = 2425153.024
$data[24,0] = 0
$data[24,1] = 0
$data[24,2] = 0
$data[25,0] = 4509324.8
$data[25,1] = 4565757.952
$data[25,2] = 4330566.144
$data[26,0] = 1303623.04
$data[26,1] = 1345964.032
$data[26,2] = 1435753.984
$data[27,0] = 49507
$data[27,1] = 43534
$data[27,2] = 42271
$data[28,0] = 481495.008
$data[28,1] = 367992
$data[28,2] = 381284.992
$data[29,0] = 25916
$data[29,1] = 32167
$data[29,2] = 25441
$data[30,0] = 579705.9840000001
$data[30,1] = 647990.976
$data[30,2] = 613358.976
$data[31,0] = 0
$data[31,1] = 0
$data[31,2] = 0
$data[32,0] = 0
$data[32,1] = 0
$data[32,2] = 0
$data[33,0] = 166999.008
$data[33,1] = 254280
$data[33,2] = 373398.016
$data[34,0] = 0
$data[34,1] = 0
$data[34,2] = 0
$data[35,0] = 0
$data[35,1] = 0
$data[35,2] = 0
$data[36,0] = 1920168.96
$data[36,1] = 1914478.976
$data[36,2] = 1571394.048
$data[37,0] = 1657506.944
$data[37,1] = 1654082.944
$data[37,2] = 1314610.048
$data[38,0] = 0
$data[38,1] = 0
$data[38,2] = 0
$data[39,0] = 108802
$data[39,1] = 97613
$data[39,2] = 98822
$data[40,0] = 67825
$data[40,1] = 70329
$data[40,2] = 70419
$data[41,0] = 0
$data[41,1] = 0
$data[41,2] = 0
$data[42,0] = 86035
$data[42,1] = 92454
$data[42,2] = 87543
$data[43,0] = 0
$data[43,1] = 0
$data[43,2] = 0
$data[44,0] = 0
$data[44,1] = 0
$data[44,2] = 0
$data[45,0] = 15319
$data[45,1] = 16796
$data[45,2] = 17866
$data[46,0] = 1270214.056
$data[46,1] = 1288518.944
$data[46,2] = 1305551.984
$data[47,0] = 875574.976
$data[47,1] = 875574.976
$data[47,2] = 875574.976
$data[48,0] = 18432
$data[48,1] = 19858
$data[48,2] = 23598
$data[49,0] = 0
$data[49,1] = 0
$data[49,2] = 0
$data[50,0] = 376207.008
$data[50,1] = 376207.008
$data[50,2] = 376207.008
$data[51,0] = 0
$data[51,1] = 16879
$data[51,2] = 30172
$data[52,0] = 0
$data[52,1] = 0
$data[52,2] = 0
$data[53,0] = 0
$data[53,1] = 0
$data[53,2] = 0
$data[54,0] = 0
$data[54,1] = 0
$data[54,2] = 0
$data[55,0] = 0
$data[55,1] = 0
$data[55,2] = 0
$data[56,0] = $null
$data[56,1] = $null
$data[56,2] = $null
$data[57,0] = $null
$data[57,1] = $null
$data[57,2] = $null
$data[58,0] = 414420.064
$data[58,1] = 403686.016
$data[58,2] = 398751.008
$data[59,0] = -96447.008
$data[59,1] = -75448
$data[59,2] = -69113
$data[60,0] = 317972.992
$data[60,1] = 328238.016
$data[60,2] = 329638.016
$data[61,0] = -130131.024
$data[61,1] = -111307
$data[61,2] = -120405
$data[62,0] = -101349
$data[62,1] = -86169
$data[62,2] = -88447
$data[63,0] = -63602
$data[63,1] = -29786
$data[63,2] = -31126
$data[64,0] = 15587
$data[64,1] = 0
$data[64,2] = 0
$data[65,0] = -65768
$data[65,1] = -27059
$data[65,2] = -39014
$data[66,0] = -341
$data[66,1] = 0
$data[66,2] = 0
$data[67,0] = -51414
$data[67,1] = -43800
$data[67,2] = -43009
$data[68,0] = 40215.008
$data[68,1] = 39412
$data[68,2] = 44058
$data[69,0] = -91629.008
$data[69,1] = -83212
$data[69,2] = -87067
$data[70,0] = $null
$data[70,1] = $null
$data[70,2] = $null
$data[71,0] = $null
$data[71,1] = $null
$data[71,2] = $null
$data[72,0] = $null
$data[72,1] = $null
$data[72,2] = $null
$data[73,0] = -79045
$data[73,1] = 30117
$data[73,2] = 7637
$data[74,0] = 26494
$data[74,1] = -9148
$data[74,2] = -6898
$data[75,0] = -4717
$data[75,1] = -2538
$data[75,2] = 14370
$data[76,0] = $null
$data[76,1] = $null
$data[76,2] = $null
$data[77,0] = $null
$data[77,1] = $null
$data[77,2] = $null
$data[78,0] = -1400
$data[78,1] = -1552
$data[78,2] = -1816
$data[79,0] = -122213
$data[79,1] = 16879
$data[79,2] = 13293

$ws.Range("BB1:BD80").Value2 = $data

Write-Output "Done writing BB1:BD80"
